$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 568.375
$ws.Range("I9").Value = 585.2857
$ws.Range("K9").Value = 585.2857
$ws.Range("M9").Value = -416.2857

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H137").Value = 2954.8
$ws.Range("I137").Value = 2073.0667
$ws.Range("J137").Value = 5600
$ws.Range("K137").Value = 6219.2001
$ws.Range("L137").Value = 16800
$ws.Range("M137").Value = -3669.2001
$ws.Range("N137").Value = -21900

$ws.Range("H138").Value = 4300.4707
$ws.Range("I138").Value = 3826.3572
$ws.Range("J138").Value = 4479.8647
$ws.Range("K138").Value = 11479.0716
$ws.Range("L138").Value = 13439.5941
$ws.Range("M138").Value = -6339.071599999999
$ws.Range("N138").Value = -23719.5941

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 220.33333
$ws.Range("I4").Value = 220.33333
$ws.Range("K4").Value = 220.33333
$ws.Range("M4").Value = -104.33333

$ws.Range("H6").Value = 254250
$ws.Range("I6").Value = 5666.6665
$ws.Range("J6").Value = 1000000
$ws.Range("K6").Value = 5666.6665
$ws.Range("L6").Value = 1000000
$ws.Range("M6").Value = -5493.6665
$ws.Range("N6").Value = -1000346

$ws.Range("H31").Value = 12752.125
$ws.Range("I31").Value = 12238.143
$ws.Range("J31").Value = 16350
$ws.Range("K31").Value = 12238.143
$ws.Range("L31").Value = 16350
$ws.Range("M31").Value = -11944.143
$ws.Range("N31").Value = -16938

$ws.Range("H32").Value = 1865029.9
$ws.Range("I32").Value = 2173021.2
$ws.Range("K32").Value = 2173021.2
$ws.Range("M32").Value = -2172734.2

$ws.Range("H61").Value = 7078.174
$ws.Range("I61").Value = 3737.6875
$ws.Range("K61").Value = 3737.6875
$ws.Range("M61").Value = -3525.6875

$ws.Range("H74").Value = 4811503.5
$ws.Range("J74").Value = 5961.0835
$ws.Range("L74").Value = 5961.0835
$ws.Range("N74").Value = -7709.0835

$ws.Range("H77").Value = 4811503.5
$ws.Range("J77").Value = 5961.0835
$ws.Range("L77").Value = 29805.4175
$ws.Range("N77").Value = -38541.4175

$ws.Range("H97").Value = 784.3
$ws.Range("I97").Value = 720.3158
$ws.Range("K97").Value = 720.3158
$ws.Range("M97").Value = -224.3158

$ws.Range("H132").Value = 373900.44
$ws.Range("I132").Value = 508270.5
$ws.Range("K132").Value = 1524811.5
$ws.Range("M132").Value = -1522281.5

$ws.Range("H136").Value = 7078.174
$ws.Range("I136").Value = 3737.6875
$ws.Range("K136").Value = 11213.0625
$ws.Range("M136").Value = -8663.0625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 664375.5600000001
$ws.Range("I134").Value = 795714.3
$ws.Range("K134").Value = 2387142.9
$ws.Range("M134").Value = -2384607.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10176
$ws.Range("I31").Value = 3775.375
$ws.Range("J31").Value = 15865.444
$ws.Range("K31").Value = 3775.375
$ws.Range("L31").Value = 15865.444
$ws.Range("M31").Value = -3480.375
$ws.Range("N31").Value = -16455.444

$ws.Range("H34").Value = 10176
$ws.Range("I34").Value = 3775.375
$ws.Range("J34").Value = 15865.444
$ws.Range("K34").Value = 3775.375
$ws.Range("L34").Value = 15865.444
$ws.Range("M34").Value = -3573.375
$ws.Range("N34").Value = -16269.444

$ws.Range("H94").Value = 1317.5555
$ws.Range("I94").Value = 800.3077
$ws.Range("K94").Value = 800.3077
$ws.Range("M94").Value = -349.3077

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 172442.48
$ws.Range("J2").Value = 31.083334
$ws.Range("L2").Value = 186.500004
$ws.Range("N2").Value = -412.500004

$ws.Range("H6").Value = 58.666668
$ws.Range("J6").Value = 199
$ws.Range("L6").Value = 597
$ws.Range("N6").Value = -823

$ws.Range("H7").Value = 3425260.8
$ws.Range("J7").Value = 3517026.8
$ws.Range("L7").Value = 10551080.4
$ws.Range("N7").Value = -10551304.4

$ws.Range("H9").Value = 593.3333
$ws.Range("I9").Value = 800
$ws.Range("J9").Value = 428
$ws.Range("K9").Value = 2400
$ws.Range("L9").Value = 1284
$ws.Range("M9").Value = -2176
$ws.Range("N9").Value = -1732

$ws.Range("H10").Value = 709.5714
$ws.Range("I10").Value = 409.4
$ws.Range("J10").Value = 1460
$ws.Range("K10").Value = 1228.2
$ws.Range("L10").Value = 4380
$ws.Range("M10").Value = -1089.2
$ws.Range("N10").Value = -4658

$ws.Range("H11").Value = 1242.5834
$ws.Range("I11").Value = 403.5
$ws.Range("J11").Value = 1410.4
$ws.Range("K11").Value = 1210.5
$ws.Range("L11").Value = 4231.200000000001
$ws.Range("M11").Value = -1070.5
$ws.Range("N11").Value = -4511.200000000001

$ws.Range("H13").Value = 1115.3
$ws.Range("I13").Value = 87
$ws.Range("J13").Value = 1556
$ws.Range("K13").Value = 261
$ws.Range("L13").Value = 4668
$ws.Range("M13").Value = -93
$ws.Range("N13").Value = -5004

$ws.Range("H38").Value = 50
$ws.Range("I38").Value = 40
$ws.Range("K38").Value = 120
$ws.Range("M38").Value = 227

$ws.Range("H46").Value = 2879.8
$ws.Range("I46").Value = 2349.75
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 7049.25
$ws.Range("L46").Value = 15000
$ws.Range("M46").Value = -6958.25
$ws.Range("N46").Value = -15182

$ws.Range("H122").Value = 1337.1072
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1337.1072
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 12033.9648
$ws.Range("N122").Value = -16933.9648
$ws.Range("M122").ClearContents()

$ws.Range("H123").Value = 5534.778
$ws.Range("I123").Value = 3380
$ws.Range("J123").Value = 9844.333000000001
$ws.Range("K123").Value = 10140
$ws.Range("L123").Value = 29532.999
$ws.Range("M123").Value = -7690
$ws.Range("N123").Value = -34432.999

$ws.Range("H131").Value = 12156
$ws.Range("I131").Value = 3011.6
$ws.Range("J131").Value = 14696.111
$ws.Range("K131").Value = 9034.799999999999
$ws.Range("L131").Value = 44088.333
$ws.Range("M131").Value = -3994.799999999999
$ws.Range("N131").Value = -54168.333

$ws.Range("H132").Value = 3552.44
$ws.Range("I132").Value = 2149.3635
$ws.Range("J132").Value = 4654.857
$ws.Range("K132").Value = 19344.2715
$ws.Range("L132").Value = 41893.713
$ws.Range("M132").Value = -16814.2715
$ws.Range("N132").Value = -46953.713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 45332
$ws.Range("J32").Value = 45332
$ws.Range("L32").Value = 45332
$ws.Range("N32").Value = -45924

$ws.Range("H70").Value = 9328.885
$ws.Range("I70").Value = 11860.667
$ws.Range("K70").Value = 11860.667
$ws.Range("M70").Value = -11590.667

$ws.Range("H73").Value = 9328.885
$ws.Range("I73").Value = 11860.667
$ws.Range("K73").Value = 11860.667
$ws.Range("M73").Value = -10924.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3726.4546
$ws.Range("I40").Value = 2999.4285
$ws.Range("J40").Value = 4998.75
$ws.Range("K40").Value = 2999.4285
$ws.Range("L40").Value = 4998.75
$ws.Range("M40").Value = -2863.4285
$ws.Range("N40").Value = -5270.75

$ws.Range("H122").Value = 47666.832
$ws.Range("J122").Value = 152627
$ws.Range("L122").Value = 457881
$ws.Range("N122").Value = -462781

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 50070

$ws.Range("H73").Value = 50070

$ws.Range("H81").Value = 2741.6667
$ws.Range("I81").Value = 2741.6667
$ws.Range("K81").Value = 5483.3334
$ws.Range("M81").Value = -4422.3334

$ws.Range("H84").Value = 2741.6667
$ws.Range("I84").Value = 2741.6667
$ws.Range("K84").Value = 27416.667
$ws.Range("M84").Value = -22112.667
